$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2794.2334
$ws.Range("J17").Value = 2382.1785
$ws.Range("L17").Value = 7146.5355
$ws.Range("N17").Value = -7482.5355
$ws.Range("H37").Value = 700
$ws.Range("I37").Value = 700
$ws.Range("K37").Value = 2100
$ws.Range("M37").Value = -1974
$ws.Range("H41").Value = 1188.7778
$ws.Range("I41").Value = 900
$ws.Range("J41").Value = 1271.2858
$ws.Range("K41").Value = 900
$ws.Range("L41").Value = 1271.2858
$ws.Range("M41").Value = -460
$ws.Range("N41").Value = -2151.2858
$ws.Range("H43").Value = 526
$ws.Range("I43").Value = 515.6
$ws.Range("J43").Value = 543.3333
$ws.Range("K43").Value = 515.6
$ws.Range("L43").Value = 543.3333
$ws.Range("M43").Value = -446.6
$ws.Range("N43").Value = -681.3333
$ws.Range("H47").Value = 20067
$ws.Range("I47").Value = 20067
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 20067
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -19095
$ws.Range("N47").ClearContents()
$ws.Range("H100").Value = 2137.9
$ws.Range("I100").Value = 914.2857
$ws.Range("J100").Value = 4993
$ws.Range("K100").Value = 914.2857
$ws.Range("L100").Value = 4993
$ws.Range("M100").Value = -373.2857
$ws.Range("N100").Value = -6075
$ws.Range("H112").Value = 1341.7428
$ws.Range("I112").Value = 500
$ws.Range("K112").Value = 1500
$ws.Range("M112").Value = -392
$ws.Range("H132").Value = 1345.6538
$ws.Range("I132").Value = 1329.9474
$ws.Range("K132").Value = 3989.8422
$ws.Range("M132").Value = -1459.8422
$ws.Range("H135").Value = 50001056
$ws.Range("I135").Value = 1081.3572
$ws.Range("K135").Value = 9732.2148
$ws.Range("M135").Value = -7197.2148
$ws.Range("H138").Value = 2176.15
$ws.Range("I138").Value = 1972.8889
$ws.Range("J138").Value = 2481.0417
$ws.Range("K138").Value = 5918.6667
$ws.Range("L138").Value = 7443.125100000001
$ws.Range("M138").Value = -778.6666999999998
$ws.Range("N138").Value = -17723.1251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 4700
$ws.Range("I31").Value = 4700
$ws.Range("K31").Value = 4700
$ws.Range("M31").Value = -4406
$ws.Range("H74").Value = 1013.35895
$ws.Range("I74").Value = 580.3
$ws.Range("J74").Value = 2456.889
$ws.Range("K74").Value = 580.3
$ws.Range("L74").Value = 2456.889
$ws.Range("M74").Value = 293.7
$ws.Range("N74").Value = -4204.889
$ws.Range("H77").Value = 1013.35895
$ws.Range("I77").Value = 580.3
$ws.Range("J77").Value = 2456.889
$ws.Range("K77").Value = 2901.5
$ws.Range("L77").Value = 12284.445
$ws.Range("M77").Value = 1466.5
$ws.Range("N77").Value = -21020.445
$ws.Range("H123").Value = 82000
$ws.Range("J123").Value = 82000
$ws.Range("L123").Value = 82000
$ws.Range("N123").Value = -91800
$ws.Range("H132").Value = 1338
$ws.Range("I132").Value = 1241.5
$ws.Range("K132").Value = 3724.5
$ws.Range("M132").Value = -1194.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H99").Value = 1060.1
$ws.Range("J99").Value = 1126
$ws.Range("L99").Value = 1126
$ws.Range("N99").Value = -4122
$ws.Range("H107").Value = 1133
$ws.Range("I107").Value = 1199
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 1199
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = 721
$ws.Range("N107").Value = -4940
$ws.Range("H134").Value = 6846.316
$ws.Range("I134").Value = 8028.933
$ws.Range("K134").Value = 24086.799
$ws.Range("M134").Value = -21551.799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 100801.664
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 100801.664
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = -100514.664
$ws.Range("N16").Value = -1374
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H93").Value = 15851.75
$ws.Range("I93").Value = 14469
$ws.Range("K93").Value = 14469
$ws.Range("M93").Value = -12597
$ws.Range("H105").Value = 1056.25
$ws.Range("I105").Value = 1056.3
$ws.Range("K105").Value = 1056.3
$ws.Range("M105").Value = 690.7
$ws.Range("I113").Value = 100801.664
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 100801.664
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = -98631.664
$ws.Range("N113").Value = -5140
$ws.Range("H132").Value = 1796.5294
$ws.Range("I132").Value = 1339.7916
$ws.Range("J132").Value = 2892.7
$ws.Range("K132").Value = 4019.3748
$ws.Range("L132").Value = 8678.099999999999
$ws.Range("M132").Value = -1489.3748
$ws.Range("N132").Value = -13738.1
$ws.Range("H134").Value = 2903.6667
$ws.Range("I134").Value = 2350.4614
$ws.Range("K134").Value = 7051.3842
$ws.Range("M134").Value = -4516.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 206.27272
$ws.Range("I2").Value = 138.88889
$ws.Range("J2").Value = 509.5
$ws.Range("K2").Value = 833.33334
$ws.Range("L2").Value = 3057
$ws.Range("M2").Value = -720.33334
$ws.Range("N2").Value = -3283
$ws.Range("H52").Value = 1500
$ws.Range("J52").Value = 1500
$ws.Range("L52").Value = 4500
$ws.Range("N52").Value = -5032
$ws.Range("H98").Value = 1304.8334
$ws.Range("I98").Value = 900
$ws.Range("J98").Value = 1385.8
$ws.Range("K98").Value = 2700
$ws.Range("L98").Value = 4157.4
$ws.Range("M98").Value = -1202
$ws.Range("N98").Value = -7153.4
$ws.Range("H107").Value = 554.3333
$ws.Range("I107").Value = 466.66666
$ws.Range("J107").Value = 583.55554
$ws.Range("K107").Value = 1399.99998
$ws.Range("L107").Value = 1750.66662
$ws.Range("M107").Value = 520.0000199999999
$ws.Range("N107").Value = -5590.66662
$ws.Range("H131").Value = 19009.64
$ws.Range("I131").Value = 825
$ws.Range("J131").Value = 19992.594
$ws.Range("K131").Value = 2475
$ws.Range("L131").Value = 59977.78200000001
$ws.Range("M131").Value = 2565
$ws.Range("N131").Value = -70057.78200000001
$ws.Range("H136").Value = 1153.2858
$ws.Range("I136").Value = 1153.2858
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3459.8574
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 1640.1426
$ws.Range("N136").ClearContents()
$ws.Range("H138").Value = 4011.182
$ws.Range("J138").Value = 4747.5713
$ws.Range("L138").Value = 14242.7139
$ws.Range("N138").Value = -24522.7139
$ws.Range("H139").Value = 1907.5
$ws.Range("I139").Value = 1904.3125
$ws.Range("J139").Value = 1933
$ws.Range("K139").Value = 5712.9375
$ws.Range("L139").Value = 5799
$ws.Range("M139").Value = -572.9375
$ws.Range("N139").Value = -16079
$ws.Range("H140").Value = 2608.4119
$ws.Range("I140").Value = 1627.7778
$ws.Range("J140").Value = 3711.625
$ws.Range("K140").Value = 4883.3334
$ws.Range("L140").Value = 11134.875
$ws.Range("M140").Value = 296.6665999999996
$ws.Range("N140").Value = -21494.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 687
$ws.Range("I80").Value = 616.3333
$ws.Range("J80").Value = 899
$ws.Range("K80").Value = 616.3333
$ws.Range("L80").Value = 899
$ws.Range("M80").Value = 381.6667
$ws.Range("N80").Value = -2895
$ws.Range("H83").Value = 687
$ws.Range("I83").Value = 616.3333
$ws.Range("J83").Value = 899
$ws.Range("K83").Value = 3081.6665
$ws.Range("L83").Value = 4495
$ws.Range("M83").Value = 1910.3335
$ws.Range("N83").Value = -14479
$ws.Range("H99").Value = 9000
$ws.Range("I99").Value = 9000
$ws.Range("K99").Value = 9000
$ws.Range("M99").Value = -6754
$ws.Range("H132").Value = 2266532.8
$ws.Range("I132").Value = 2962135.5
$ws.Range("K132").Value = 8886406.5
$ws.Range("M132").Value = -8883876.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 10000
$ws.Range("I30").Value = 10000
$ws.Range("K30").Value = 10000
$ws.Range("M30").Value = -9892
$ws.Range("H46").Value = 1547.2222
$ws.Range("I46").Value = 1085.4445
$ws.Range("J46").Value = 2009
$ws.Range("K46").Value = 1085.4445
$ws.Range("L46").Value = 2009
$ws.Range("M46").Value = -897.4445000000001
$ws.Range("N46").Value = -2385
$ws.Range("H132").Value = 2350.6956
$ws.Range("I132").Value = 1601.1666
$ws.Range("J132").Value = 2615.2354
$ws.Range("K132").Value = 4803.4998
$ws.Range("L132").Value = 7845.706200000001
$ws.Range("M132").Value = -2273.4998
$ws.Range("N132").Value = -12905.7062
$ws.Range("H136").Value = 4212.0454
$ws.Range("I136").Value = 3408.9473
$ws.Range("J136").Value = 9298.333000000001
$ws.Range("K136").Value = 10226.8419
$ws.Range("L136").Value = 27894.999
$ws.Range("M136").Value = -7676.841899999999
$ws.Range("N136").Value = -32994.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1763.5
$ws.Range("I96").Value = 990
$ws.Range("K96").Value = 990
$ws.Range("M96").Value = 383
$ws.Range("H100").Value = 1075.1666
$ws.Range("I100").Value = 930.2
$ws.Range("K100").Value = 1860.4
$ws.Range("M100").Value = -1319.4
$ws.Range("H132").Value = 1570.091
$ws.Range("I132").Value = 1034.8518
$ws.Range("J132").Value = 3978.6667
$ws.Range("K132").Value = 3104.5554
$ws.Range("L132").Value = 11936.0001
$ws.Range("M132").Value = -574.5553999999997
$ws.Range("N132").Value = -16996.0001
$ws.Range("H136").Value = 15434000
$ws.Range("I136").Value = 25254404
$ws.Range("K136").Value = 75763212
$ws.Range("M136").Value = -75760662
